# ----------------------------------------------------------------------
# updated batch additional test scenarios
#
# - updateByBatchId: selection moves from C15 to A20 (and loses tab focus)
# - new worksheet "deleteBatchByBatchId" inserted before "getBatchById",
#   populated with new test-scenario rows, becomes the active/selected tab
# - getBatchById: selection moves from D13 to C3
# - getBatchByName: selection moves from F5 to E4
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- updateByBatchId: change selection, no longer the active tab ---
$wsUpdate = $wb.Worksheets.Item("updateByBatchId")
$wsUpdate.Select()
$wsUpdate.Range("A20").Select()

# --- getBatchById: change selection ---
$wsById = $wb.Worksheets.Item("getBatchById")
$wsById.Range("C3").Select()

# --- getBatchByName: change selection ---
$wsByName = $wb.Worksheets.Item("getBatchByName")
$wsByName.Range("E4").Select()

# --- insert the new "deleteBatchByBatchId" sheet right before getBatchById ---
$newSheet = $wb.Worksheets.Add($wsById)
$newSheet.Name = "deleteBatchByBatchId"

$newSheet.Range("A1").Value = "baseUrl"
$newSheet.Range("B1").Value = "endPoint"
$newSheet.Range("C1").Value = "AuthStatus"
$newSheet.Range("D1").Value = "httpStatusCode"
$newSheet.Range("E1").Value = "Scenario"
$newSheet.Range("F1").Value = "BatchId"

$newSheet.Range("A2").Value = "https://lms-marchapi-hackathon-a258d2bbd43b.herokuapp.com/lms"
$newSheet.Range("B2").Value = "/batches/"
$newSheet.Range("C2").Value = "withBasicAuth"
$newSheet.Range("D2").Value = "200-OK"
$newSheet.Range("E2").Value = "positive"
$newSheet.Range("F2").Value = 8898

$newSheet.Range("A3").Value = "https://lms-marchapi-hackathon-a258d2bbd43b.herokuapp.com/lms"
$newSheet.Range("B3").Value = "/batches3344/"
$newSheet.Range("C3").Value = "withBasicAuth"
$newSheet.Range("D3").Value = "404-  Not Found"
$newSheet.Range("E3").Value = "negative"
$newSheet.Range("F3").Value = 9018

$newSheet.Range("A4").Value = "https://lms-marchapi-hackathon-a258d2bbd43b.herokuapp.com/lms"
$newSheet.Range("B4").Value = "/batches/"
$newSheet.Range("C4").Value = "withBasicAuth"
$newSheet.Range("D4").Value = "404-  Not Found"
$newSheet.Range("E4").Value = "negative"
$newSheet.Range("F4").Value = 8954

$newSheet.Range("A5").Value = "https://lms-marchapi-hackathon-a258d2bbd43b.herokuapp.com/lms"
$newSheet.Range("B5").Value = "/batches/"
$newSheet.Range("C5").Value = "NoAuth"
$newSheet.Range("D5").Value = "401-Unauthorized"
$newSheet.Range("E5").Value = "negative"
$newSheet.Range("F5").Value = 8876

# whole table is formatted as Text, matching the original workbook's style
$newSheet.Range("A1:F5").NumberFormat = "@"

# best-effort column widths (engine quantizes to 1/6 character, so these
# are the closest achievable approximations of the authored best-fit widths)
$newSheet.Columns.Item(1).ColumnWidth = 58.83333333333333
$newSheet.Columns.Item(2).ColumnWidth = 15.333333333333334
$newSheet.Columns.Item(3).ColumnWidth = 9.0
$newSheet.Columns.Item(4).ColumnWidth = 15.0
$newSheet.Columns.Item(5).ColumnWidth = 7.0
$newSheet.Columns.Item(6).ColumnWidth = 6.5

# the new sheet becomes the active / selected tab, with column F fully selected
$newSheet.Select()
$newSheet.Range("F1:F1048576").Select()
